$d = $word.ActiveDocument

$replacements = @(
    @("994×4=3976", "605×7=4235"),
    @("900×6=5400", "358×8=2864"),
    @("310×6=1860", "408×3=1224"),
    @("636×2=1272", "150×5=750"),
    @("802×6=4812", "232×8=1856"),
    @("994×7=6958", "625×8=5000"),
    @("363×4=1452", "539×5=2695"),
    @("912×3=2736", "676×3=2028"),
    @("392×6=2352", "811×8=6488"),
    @("438×4=1752", "715×8=5720"),
    @("243×2=486",  "340×6=2040"),
    @("764×5=3820", "902×9=8118"),
    @("861×9=7749", "722×2=1444"),
    @("215×5=1075", "521×4=2084"),
    @("285×3=855",  "903×5=4515"),
    @("317×5=1585", "137×2=274"),
    @("835×2=1670", "148×2=296"),
    @("957×5=4785", "729×6=4374"),
    @("610×3=1830", "748×9=6732"),
    @("594×4=2376", "232×3=696"),
    @("357×9=3213", "449×7=3143"),
    @("336×3=1008", "162×4=648"),
    @("603×8=4824", "722×6=4332"),
    @("855×2=1710", "803×9=7227"),
    @("396×6=2376", "238×6=1428")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
